$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear the "Kernel time (ms)" raw data in column B (rows 3-5), keep header B2
$ws.Range("B3:B5").ClearContents()

# 2. Clear K2:L10 entirely (values+format; nothing worth keeping there)
$ws.Range("K2:L10").Clear()

# 3. Remove the yellow fill highlighting: first on K1 alone (mutates its style in place),
#    then apply the same "no fill" to A3:A5 so they consolidate onto the very same style index.
$ws.Range("K1").Interior.ColorIndex = -4142
$ws.Range("A3:A5").Interior.ColorIndex = -4142

# 4. Clear K1's value/text (keep the now-updated format)
$ws.Range("K1").ClearContents()

# 5. Update H column with the new benchmark values
$ws.Range("H3").Value = 1768874.44
$ws.Range("H4").Value = 221633.898
$ws.Range("H5").Value = 27881.647
$ws.Range("H6").Value = 3489.747
$ws.Range("H8").Value = 3516.354
$ws.Range("H9").Value = 27818.054
$ws.Range("H10").Value = 221552.069

# 6. Update sheet view: clear frozen/top-left cell and change selection
$ws.Range("D8").Select()
